$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at C (old Latitude/Longitude shift right to D/E).
$ws.Columns.Item(3).Insert()

# Header row: B1 becomes "AudioStatus" (new data), C1 becomes "WifiStatus"
# (holds what used to be the "Status" column values), D1/E1 keep Latitude/Longitude.
$ws.Cells.Item(1, 2).Value = "AudioStatus"
$ws.Cells.Item(1, 3).Value = "WifiStatus"

# Per-building AudioStatus and WifiStatus values.
# WifiStatus simply mirrors the previous "Status" column (now in column B).
$buildings = @(
    "Holden Hall",
    "West Hall",
    "Admin Building",
    "Student Union Building",
    "Library",
    "Agricultural Pavilion",
    "Special Collections Library",
    "Food Technology Building",
    "Psychology Building",
    "Gates Hall",
    "Wall Hall",
    "Knapp Hall",
    "Learning Hall",
    "Food Pantry",
    "McClellan Hall",
    "Sneed Hall",
    "Livermore Center",
    "Experimental Sciences Building",
    "Madoxx Engineering Research Center"
)

$audioStatus = @(
    "Inactive",
    "Inactive",
    "Active",
    "Active",
    "Inactive",
    "Active",
    "Active",
    "Inactive",
    "Active",
    "Inactive",
    "Active",
    "Inactive",
    "Active",
    "Inactive",
    "Active",
    "Active",
    "Active",
    "Active",
    "Inactive"
)

$wifiStatus = @(
    "Active",
    "Active",
    "Inactive",
    "Inactive",
    "Active",
    "Active",
    "Inactive",
    "Active",
    "Active",
    "Inactive",
    "Active",
    "Inactive",
    "Active",
    "Inactive",
    "Active",
    "Active",
    "Active",
    "Active",
    "Inactive"
)

for ($i = 0; $i -lt $buildings.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $audioStatus[$i]
    $ws.Cells.Item($row, 3).Value = $wifiStatus[$i]
}

# All rows (including the header) now share the same row height.
$ws.Range("A1:E20").RowHeight = 19.5

# Column widths/formatting for the new column should match the other text columns.
$ws.Columns.Item(3).ColumnWidth = $ws.Columns.Item(2).ColumnWidth
